{"js": "// Update the title/date paragraph (first paragraph in the document body).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titleParagraph = paragraphs.items[0];\ntitleParagraph.insertText(\"2025-08-01 Friday\", Word.InsertLocation.replace);\n\n// Update every equation cell in the first (only) table at once, in\n// document (row-major) order -- one array entry per row, left to right --\n// which preserves each cell's existing run formatting (font/size) because\n// only the text runs are rewritten, not the cell/paragraph structure.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst newValues = [\n  [\"91-76=\", \"49-15=\", \"85-74=\", \"46-3=\", \"13+22=\"],\n  [\"18+40=\", \"59-4=\", \"6+29=\", \"47-38=\", \"87-30=\"],\n  [\"89-67=\", \"43+32=\", \"6+85=\", \"98-98=\", \"88-54=\"],\n  [\"61-48=\", \"71-33=\", \"15+77=\", \"18+20=\", \"29+26=\"],\n  [\"83-38=\", \"85-21=\", \"50+42=\", \"28+25=\", \"33+37=\"],\n  [\"23-15=\", \"78+8=\", \"45-18=\", \"55-29=\", \"70+5=\"],\n  [\"97-28=\", \"46-17=\", \"24+67=\", \"52+35=\", \"99-27=\"],\n  [\"47-23=\", \"6-6=\", \"32-28=\", \"20+15=\", \"69-28=\"],\n  [\"6+64=\", \"83-39=\", \"6+13=\", \"63-60=\", \"64-30=\"],\n  [\"47-33=\", \"7+52=\", \"33-20=\", \"26-10=\", \"2+86=\"],\n  [\"68-30=\", \"61+29=\", \"27-2=\", \"54+11=\", \"75-48=\"],\n  [\"57-26=\", \"14-3=\", \"32-8=\", \"65-60=\", \"10-8=\"],\n  [\"16-0=\", \"83-82=\", \"96-92=\", \"69-52=\", \"45+20=\"],\n  [\"95-19=\", \"91-64=\", \"5+56=\", \"87+5=\", \"85-72=\"],\n  [\"54-32=\", \"26+61=\", \"7+63=\", \"49-20=\", \"14+82=\"],\n  [\"2+0=\", \"91-76=\", \"39-20=\", \"59-58=\", \"91-52=\"],\n  [\"9+59=\", \"89-57=\", \"35+44=\", \"82+4=\", \"20+12=\"],\n  [\"96-28=\", \"45-30=\", \"68+10=\", \"17+20=\", \"21-19=\"],\n  [\"94-89=\", \"54-34=\", \"11+36=\", \"83-35=\", \"33+0=\"],\n  [\"68-30=\", \"4+5=\", \"93-90=\", \"89-86=\", \"8+26=\"]\n];\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the title/date paragraph (first paragraph in the document body).\n$titleOld = \"2025-07-31 Thursday\"\n$titleNew = \"2025-08-01 Friday\"\n$p1 = $d.Paragraphs.Item(1)\nif ($p1.Range.Text.TrimEnd([char]13, [char]7) -eq $titleOld) {\n    $p1.Range.Text = $titleNew\n}\n\n# Update every equation cell in the first table, in document (row-major) order.\n# Each entry is the expected old text followed by its replacement; this\n# preserves correctness even though some \"old\" values repeat with\n# different replacements at different positions.\n$pairs = @(\n    @(\"13+6=\",\"91-76=\"),\n    @(\"80+17=\",\"49-15=\"),\n    @(\"19+10=\",\"85-74=\"),\n    @(\"89-82=\",\"46-3=\"),\n    @(\"10+31=\",\"13+22=\"),\n    @(\"17+11=\",\"18+40=\"),\n    @(\"89-68=\",\"59-4=\"),\n    @(\"50-19=\",\"6+29=\"),\n    @(\"41-4=\",\"47-38=\"),\n    @(\"71-1=\",\"87-30=\"),\n    @(\"95-68=\",\"89-67=\"),\n    @(\"61+9=\",\"43+32=\"),\n    @(\"28-7=\",\"6+85=\"),\n    @(\"40+13=\",\"98-98=\"),\n    @(\"8-8=\",\"88-54=\"),\n    @(\"46-19=\",\"61-48=\"),\n    @(\"36-6=\",\"71-33=\"),\n    @(\"43+30=\",\"15+77=\"),\n    @(\"18+76=\",\"18+20=\"),\n    @(\"36+37=\",\"29+26=\"),\n    @(\"21+33=\",\"83-38=\"),\n    @(\"50-34=\",\"85-21=\"),\n    @(\"27+36=\",\"50+42=\"),\n    @(\"81-75=\",\"28+25=\"),\n    @(\"63+12=\",\"33+37=\"),\n    @(\"88-34=\",\"23-15=\"),\n    @(\"84-83=\",\"78+8=\"),\n    @(\"83-67=\",\"45-18=\"),\n    @(\"30+44=\",\"55-29=\"),\n    @(\"35-19=\",\"70+5=\"),\n    @(\"72-25=\",\"97-28=\"),\n    @(\"28+2=\",\"46-17=\"),\n    @(\"97-31=\",\"24+67=\"),\n    @(\"40-3=\",\"52+35=\"),\n    @(\"73-2=\",\"99-27=\"),\n    @(\"75-66=\",\"47-23=\"),\n    @(\"34+27=\",\"6-6=\"),\n    @(\"17+77=\",\"32-28=\"),\n    @(\"93-17=\",\"20+15=\"),\n    @(\"85-75=\",\"69-28=\"),\n    @(\"60-25=\",\"6+64=\"),\n    @(\"10+57=\",\"83-39=\"),\n    @(\"90-57=\",\"6+13=\"),\n    @(\"33+22=\",\"63-60=\"),\n    @(\"24-4=\",\"64-30=\"),\n    @(\"43-40=\",\"47-33=\"),\n    @(\"61-7=\",\"7+52=\"),\n    @(\"53+21=\",\"33-20=\"),\n    @(\"93-88=\",\"26-10=\"),\n    @(\"95-86=\",\"2+86=\"),\n    @(\"63-12=\",\"68-30=\"),\n    @(\"36+55=\",\"61+29=\"),\n    @(\"50-26=\",\"27-2=\"),\n    @(\"27+27=\",\"54+11=\"),\n    @(\"35+48=\",\"75-48=\"),\n    @(\"80-46=\",\"57-26=\"),\n    @(\"97-8=\",\"14-3=\"),\n    @(\"49-42=\",\"32-8=\"),\n    @(\"47+39=\",\"65-60=\"),\n    @(\"30-29=\",\"10-8=\"),\n    @(\"27-21=\",\"16-0=\"),\n    @(\"66-27=\",\"83-82=\"),\n    @(\"88+11=\",\"96-92=\"),\n    @(\"35+9=\",\"69-52=\"),\n    @(\"1+86=\",\"45+20=\"),\n    @(\"11+72=\",\"95-19=\"),\n    @(\"86-22=\",\"91-64=\"),\n    @(\"1+74=\",\"5+56=\"),\n    @(\"81-39=\",\"87+5=\"),\n    @(\"12+42=\",\"85-72=\"),\n    @(\"12+30=\",\"54-32=\"),\n    @(\"88-29=\",\"26+61=\"),\n    @(\"21+41=\",\"7+63=\"),\n    @(\"39+53=\",\"49-20=\"),\n    @(\"55-34=\",\"14+82=\"),\n    @(\"53-26=\",\"2+0=\"),\n    @(\"96-18=\",\"91-76=\"),\n    @(\"50-41=\",\"39-20=\"),\n    @(\"86-56=\",\"59-58=\"),\n    @(\"23+20=\",\"91-52=\"),\n    @(\"32+18=\",\"9+59=\"),\n    @(\"15+67=\",\"89-57=\"),\n    @(\"92-66=\",\"35+44=\"),\n    @(\"67-62=\",\"82+4=\"),\n    @(\"44+23=\",\"20+12=\"),\n    @(\"95-73=\",\"96-28=\"),\n    @(\"56-42=\",\"45-30=\"),\n    @(\"9+35=\",\"68+10=\"),\n    @(\"72-33=\",\"17+20=\"),\n    @(\"67-64=\",\"21-19=\"),\n    @(\"92-19=\",\"94-89=\"),\n    @(\"68-0=\",\"54-34=\"),\n    @(\"23+50=\",\"11+36=\"),\n    @(\"30-14=\",\"83-35=\"),\n    @(\"35-3=\",\"33+0=\"),\n    @(\"35+48=\",\"68-30=\"),\n    @(\"58+1=\",\"4+5=\"),\n    @(\"2+17=\",\"93-90=\"),\n    @(\"91-2=\",\"89-86=\"),\n    @(\"6+44=\",\"8+26=\")\n)\n\n$tbl = $d.Tables.Item(1)\n$rows = $tbl.Rows.Count\n$cols = $tbl.Columns.Count\n\n$i = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $cell = $tbl.Cell($r, $c)\n        $newText = $pairs[$i][1]\n        $cell.Range.Text = $newText\n        $i++\n    }\n}\n"}
